$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old layout repeated 38 numeric columns (B:AN); the new layout only needs 18 (B:T).
# Clear out the now-unused trailing columns U:AN for rows 1-19 before rewriting the data.
$ws.Range("U1:AN19").Clear()

$row1 = New-Object 'object[,]' 1,19
$row1[0,0] = 0
$row1[0,1] = 1
$row1[0,2] = 2
$row1[0,3] = 3
$row1[0,4] = 4
$row1[0,5] = 5
$row1[0,6] = 6
$row1[0,7] = 7
$row1[0,8] = 8
$row1[0,9] = 9
$row1[0,10] = 10
$row1[0,11] = 11
$row1[0,12] = 12
$row1[0,13] = 13
$row1[0,14] = 14
$row1[0,15] = 15
$row1[0,16] = 16
$row1[0,17] = 17
$row1[0,18] = 18
$ws.Range("B1:T1").Value = $row1

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = 0
$row2[0,1] = "HKL"
$row2[0,2] = "[2, 2, 0]"
$row2[0,3] = "[2, 0, 0]"
$row2[0,4] = "[2, 1, 1]"
$row2[0,5] = "[4, 0, 0]"
$row2[0,6] = "[3, 2, 1]"
$row2[0,7] = "[1, 1, 0]"
$row2[0,8] = "[2, 2, 2]"
$row2[0,9] = "[3, 1, 0]"
$row2[0,10] = "1Pair-A"
$row2[0,11] = "1Pair-B"
$row2[0,12] = "2Pairs-A"
$row2[0,13] = "2Pairs-B"
$row2[0,14] = "3Pairs-A"
$row2[0,15] = "3Pairs-B"
$row2[0,16] = "3Pairs-C"
$row2[0,17] = "4Pairs"
$row2[0,18] = "5A4F"
$row2[0,19] = "MaxUnique"
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = 1
$row3[0,1] = "BT8Hex_2.5"
$row3[0,2] = 1.000186219523844
$row3[0,3] = 0.9992551188448509
$row3[0,4] = 1.000186219523844
$row3[0,5] = 0.9992551188448509
$row3[0,6] = 1.000186219523844
$row3[0,7] = 1.000186219523844
$row3[0,8] = 1.000496588078894
$row3[0,9] = 0.9995903150010274
$row3[0,10] = 1.000186219523844
$row3[0,11] = 1.000186219523844
$row3[0,12] = 0.9997206691843473
$row3[0,13] = 0.9997206691843473
$row3[0,14] = 0.9996772177899073
$row3[0,15] = 0.9998758526308462
$row3[0,16] = 0.9998758526308462
$row3[0,17] = 0.9999534443540956
$row3[0,18] = 0.9999534443540956
$row3[0,19] = 0.999983446749384
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = 2
$row4[0,1] = "BT8Hex_5"
$row4[0,2] = 1.000359494649805
$row4[0,3] = 0.9985620152088167
$row4[0,4] = 1.000359494649805
$row4[0,5] = 0.9985620152088167
$row4[0,6] = 1.000359494649805
$row4[0,7] = 1.000359494649805
$row4[0,8] = 1.000958657193137
$row4[0,9] = 0.9992091083749851
$row4[0,10] = 1.000359494649805
$row4[0,11] = 1.000359494649805
$row4[0,12] = 0.9994607549293111
$row4[0,13] = 0.9994607549293111
$row4[0,14] = 0.9993768727445357
$row4[0,15] = 0.9997603348361425
$row4[0,16] = 0.9997603348361425
$row4[0,17] = 0.9999101247895581
$row4[0,18] = 0.9999101247895581
$row4[0,19] = 0.9999680441210591
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = 3
$row5[0,1] = "BT8Hex_10"
$row5[0,2] = 1.000688928348568
$row5[0,3] = 0.9972442770850863
$row5[0,4] = 1.000688928348568
$row5[0,5] = 0.9972442770850863
$row5[0,6] = 1.000688928348568
$row5[0,7] = 1.000688928348568
$row5[0,8] = 1.001837149048359
$row5[0,9] = 0.998484351190449
$row5[0,10] = 1.000688928348568
$row5[0,11] = 1.000688928348568
$row5[0,12] = 0.9989666027168274
$row5[0,13] = 0.9989666027168274
$row5[0,14] = 0.9988058522080346
$row5[0,15] = 0.999540711260741
$row5[0,16] = 0.999540711260741
$row5[0,17] = 0.9998277655326977
$row5[0,18] = 0.9998277655326977
$row5[0,19] = 0.999938760394933
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = 4
$row6[0,1] = "BT8Hex_15"
$row6[0,2] = 1.00101044000263
$row6[0,3] = 0.9959582394373682
$row6[0,4] = 1.00101044000263
$row6[0,5] = 0.9959582394373682
$row6[0,6] = 1.00101044000263
$row6[0,7] = 1.00101044000263
$row6[0,8] = 1.002694512966136
$row6[0,9] = 0.9977770326001525
$row6[0,10] = 1.00101044000263
$row6[0,11] = 1.00101044000263
$row6[0,12] = 0.998484339719999
$row6[0,13] = 0.998484339719999
$row6[0,14] = 0.9982485706800501
$row6[0,15] = 0.9993263731475426
$row6[0,16] = 0.9993263731475426
$row6[0,17] = 0.9997473898613144
$row6[0,18] = 0.9997473898613144
$row6[0,19] = 0.999910184168591
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = 5
$row7[0,1] = "Spiral2.5"
$row7[0,2] = 1.000013275794241
$row7[0,3] = 0.9999468950444055
$row7[0,4] = 1.000013275794241
$row7[0,5] = 0.9999468950444055
$row7[0,6] = 1.000013275794241
$row7[0,7] = 1.000013275794241
$row7[0,8] = 1.000035403618783
$row7[0,9] = 0.9999707918808989
$row7[0,10] = 1.000013275794241
$row7[0,11] = 1.000013275794241
$row7[0,12] = 0.9999800854193233
$row7[0,13] = 0.9999800854193233
$row7[0,14] = 0.9999769875731818
$row7[0,15] = 0.9999911488776293
$row7[0,16] = 0.9999911488776293
$row7[0,17] = 0.9999966806067822
$row7[0,18] = 0.9999966806067822
$row7[0,19] = 0.9999988196544685
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = 6
$row8[0,1] = "Spiral5"
$row8[0,2] = 1.000033838354057
$row8[0,3] = 0.9998646445687531
$row8[0,4] = 1.000033838354057
$row8[0,5] = 0.9998646445687531
$row8[0,6] = 1.000033838354057
$row8[0,7] = 1.000033838354057
$row8[0,8] = 1.000090236995333
$row8[0,9] = 0.9999255542012787
$row8[0,10] = 1.000033838354057
$row8[0,11] = 1.000033838354057
$row8[0,12] = 0.9999492414614053
$row8[0,13] = 0.9999492414614053
$row8[0,14] = 0.9999413457080298
$row8[0,15] = 0.9999774404256226
$row8[0,16] = 0.9999774404256226
$row8[0,17] = 0.9999915399077314
$row8[0,18] = 0.9999915399077314
$row8[0,19] = 0.9999969918045896
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = 7
$row9[0,1] = "Spiral7.5"
$row9[0,2] = 1.000047378675377
$row9[0,3] = 0.999810482714136
$row9[0,4] = 1.000047378675377
$row9[0,5] = 0.999810482714136
$row9[0,6] = 1.000047378675377
$row9[0,7] = 1.000047378675377
$row9[0,8] = 1.000126345384839
$row9[0,9] = 0.9998957656465441
$row9[0,10] = 1.000047378675377
$row9[0,11] = 1.000047378675377
$row9[0,12] = 0.9999289306947565
$row9[0,13] = 0.9999289306947565
$row9[0,14] = 0.9999178756786856
$row9[0,15] = 0.9999684133549632
$row9[0,16] = 0.9999684133549632
$row9[0,17] = 0.9999881546850666
$row9[0,18] = 0.9999881546850666
$row9[0,19] = 0.9999957882952749
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = 8
$row10[0,1] = "Spiral10"
$row10[0,2] = 1.000102161755974
$row10[0,3] = 0.9995913553296635
$row10[0,4] = 1.000102161755974
$row10[0,5] = 0.9995913553296635
$row10[0,6] = 1.000102161755974
$row10[0,7] = 1.000102161755974
$row10[0,8] = 1.000272430465849
$row10[0,9] = 0.9997752447810944
$row10[0,10] = 1.000102161755974
$row10[0,11] = 1.000102161755974
$row10[0,12] = 0.9998467585428186
$row10[0,13] = 0.9998467585428186
$row10[0,14] = 0.9998229206222439
$row10[0,15] = 0.9999318929472037
$row10[0,16] = 0.9999318929472037
$row10[0,17] = 0.9999744601493963
$row10[0,18] = 0.9999744601493963
$row10[0,19] = 0.9999909193074216
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = 9
$row11[0,1] = "Spiral15"
$row11[0,2] = 1.000169461465969
$row11[0,3] = 0.9993221533427838
$row11[0,4] = 1.000169461465969
$row11[0,5] = 0.9993221533427838
$row11[0,6] = 1.000169461465969
$row11[0,7] = 1.000169461465969
$row11[0,8] = 1.000451896201525
$row11[0,9] = 0.9996271839543713
$row11[0,10] = 1.000169461465969
$row11[0,11] = 1.000169461465969
$row11[0,12] = 0.9997458074043765
$row11[0,13] = 0.9997458074043765
$row11[0,14] = 0.9997062662543748
$row11[0,15] = 0.9998870254249074
$row11[0,16] = 0.9998870254249074
$row11[0,17] = 0.9999576344351728
$row11[0,18] = 0.9999576344351728
$row11[0,19] = 0.9999849363160979
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = 10
$row12[0,1] = "OffsetF45"
$row12[0,2] = 0.9984980511806822
$row12[0,3] = 1.006007798933863
$row12[0,4] = 0.9984980511806822
$row12[0,5] = 1.006007798933863
$row12[0,6] = 0.9984980511806822
$row12[0,7] = 0.9984980511806822
$row12[0,8] = 0.9959947961708796
$row12[0,9] = 1.003304289263805
$row12[0,10] = 0.9984980511806822
$row12[0,11] = 0.9984980511806822
$row12[0,12] = 1.002252925057273
$row12[0,13] = 1.002252925057273
$row12[0,14] = 1.002603379792784
$row12[0,15] = 1.001001300431743
$row12[0,16] = 1.001001300431743
$row12[0,17] = 1.000375488118978
$row12[0,18] = 1.000375488118978
$row12[0,19] = 1.000133506318433
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = 11
$row13[0,1] = "OffsetA45"
$row13[0,2] = 0.9996122202395928
$row13[0,3] = 1.001551139753336
$row13[0,4] = 0.9996122202395928
$row13[0,5] = 1.001551139753336
$row13[0,6] = 0.9996122202395928
$row13[0,7] = 0.9996122202395928
$row13[0,8] = 0.9989659102351678
$row13[0,9] = 1.000853126423775
$row13[0,10] = 0.9996122202395928
$row13[0,11] = 0.9996122202395928
$row13[0,12] = 1.000581679996464
$row13[0,13] = 1.000581679996464
$row13[0,14] = 1.000672162138901
$row13[0,15] = 1.000258526744174
$row13[0,16] = 1.000258526744174
$row13[0,17] = 1.000096950118029
$row13[0,18] = 1.000096950118029
$row13[0,19] = 1.000034472855176
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,20
$row14[0,0] = 12
$row14[0,1] = "OffsetFTD"
$row14[0,2] = 0.9984981356332555
$row14[0,3] = 1.006007490867104
$row14[0,4] = 0.9984981356332555
$row14[0,5] = 1.006007490867104
$row14[0,6] = 0.9984981356332555
$row14[0,7] = 0.9984981356332555
$row14[0,8] = 0.9959950108599505
$row14[0,9] = 1.00330411473018
$row14[0,10] = 0.9984981356332555
$row14[0,11] = 0.9984981356332555
$row14[0,12] = 1.00225281325018
$row14[0,13] = 1.00225281325018
$row14[0,14] = 1.002603247076847
$row14[0,15] = 1.001001254044538
$row14[0,16] = 1.001001254044538
$row14[0,17] = 1.000375474441718
$row14[0,18] = 1.000375474441718
$row14[0,19] = 1.000133503892833
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,20
$row15[0,0] = 13
$row15[0,1] = "OffsetATD"
$row15[0,2] = 0.9996122291325652
$row15[0,3] = 1.001551060553058
$row15[0,4] = 0.9996122291325652
$row15[0,5] = 1.001551060553058
$row15[0,6] = 0.9996122291325652
$row15[0,7] = 0.9996122291325652
$row15[0,8] = 0.9989659596189684
$row15[0,9] = 1.000853086972874
$row15[0,10] = 0.9996122291325652
$row15[0,11] = 0.9996122291325652
$row15[0,12] = 1.000581644842812
$row15[0,13] = 1.000581644842812
$row15[0,14] = 1.000672125552833
$row15[0,15] = 1.00025850627273
$row15[0,16] = 1.00025850627273
$row15[0,17] = 1.000096936987688
$row15[0,18] = 1.000096936987688
$row15[0,19] = 1.0000344657571
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,20
$row16[0,0] = 14
$row16[0,1] = "Holden2.5"
$row16[0,2] = 1.003632224036981
$row16[0,3] = 0.9854710894363202
$row16[0,4] = 1.003632224036981
$row16[0,5] = 0.9854710894363202
$row16[0,6] = 1.003632224036981
$row16[0,7] = 1.003632224036981
$row16[0,8] = 1.009685946805441
$row16[0,9] = 0.9920090990811403
$row16[0,10] = 1.003632224036981
$row16[0,11] = 1.003632224036981
$row16[0,12] = 0.9945516567366508
$row16[0,13] = 0.9945516567366508
$row16[0,14] = 0.9937041375181472
$row16[0,15] = 0.9975785125034277
$row16[0,16] = 0.9975785125034277
$row16[0,17] = 0.9990919403868161
$row16[0,18] = 0.9990919403868161
$row16[0,19] = 0.9996771345723077
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,20
$row17[0,0] = 15
$row17[0,1] = "Holden5"
$row17[0,2] = 1.002972234468734
$row17[0,3] = 0.9881110580767984
$row17[0,4] = 1.002972234468734
$row17[0,5] = 0.9881110580767984
$row17[0,6] = 1.002972234468734
$row17[0,7] = 1.002972234468734
$row17[0,8] = 1.007925965600576
$row17[0,9] = 0.9934610817393195
$row17[0,10] = 1.002972234468734
$row17[0,11] = 1.002972234468734
$row17[0,12] = 0.9955416462727662
$row17[0,13] = 0.9955416462727662
$row17[0,14] = 0.9948481247616172
$row17[0,15] = 0.9980185090047554
$row17[0,16] = 0.9980185090047554
$row17[0,17] = 0.9992569403707501
$row17[0,18] = 0.9992569403707501
$row17[0,19] = 0.9997358014704827
$ws.Range("A17:T17").Value = $row17

$row18 = New-Object 'object[,]' 1,20
$row18[0,0] = 16
$row18[0,1] = "Holden10"
$row18[0,2] = 1.001645132917943
$row18[0,3] = 0.9934194849595119
$row18[0,4] = 1.001645132917943
$row18[0,5] = 0.9934194849595119
$row18[0,6] = 1.001645132917943
$row18[0,7] = 1.001645132917943
$row18[0,8] = 1.004387007752333
$row18[0,9] = 0.9963807164382749
$row18[0,10] = 1.001645132917943
$row18[0,11] = 1.001645132917943
$row18[0,12] = 0.9975323089387275
$row18[0,13] = 0.9975323089387275
$row18[0,14] = 0.99714844477191
$row18[0,15] = 0.9989032502651326
$row18[0,16] = 0.9989032502651326
$row18[0,17] = 0.9995887209283352
$row18[0,18] = 0.9995887209283352
$row18[0,19] = 0.9998537679839914
$ws.Range("A18:T18").Value = $row18

$row19 = New-Object 'object[,]' 1,20
$row19[0,0] = 17
$row19[0,1] = "Holden15"
$row19[0,2] = 1.001795853925422
$row19[0,3] = 0.9928166108433267
$row19[0,4] = 1.001795853925422
$row19[0,5] = 0.9928166108433267
$row19[0,6] = 1.001795853925422
$row19[0,7] = 1.001795853925422
$row19[0,8] = 1.004788927645208
$row19[0,9] = 0.9960491322077251
$row19[0,10] = 1.001795853925422
$row19[0,11] = 1.001795853925422
$row19[0,12] = 0.9973062323843742
$row19[0,13] = 0.9973062323843742
$row19[0,14] = 0.9968871989921578
$row19[0,15] = 0.9988027728980567
$row19[0,16] = 0.9988027728980567
$row19[0,17] = 0.999551043154898
$row19[0,18] = 0.999551043154898
$row19[0,19] = 0.9998403720787542
$ws.Range("A19:T19").Value = $row19

$row20 = New-Object 'object[,]' 1,20
$row20[0,0] = 18
$row20[0,1] = "HexGrid-90degTilt2.5degRes"
$row20[0,2] = 1.000000768512052
$row20[0,3] = 0.999996924088603
$row20[0,4] = 1.000000768512052
$row20[0,5] = 0.999996924088603
$row20[0,6] = 1.000000768512052
$row20[0,7] = 1.000000768512052
$row20[0,8] = 1.000002050868922
$row20[0,9] = 0.9999983078362267
$row20[0,10] = 1.000000768512052
$row20[0,11] = 1.000000768512052
$row20[0,12] = 0.9999988463003275
$row20[0,13] = 0.9999988463003275
$row20[0,14] = 0.999998666812294
$row20[0,15] = 0.9999994870375691
$row20[0,16] = 0.9999994870375691
$row20[0,17] = 0.9999998074061899
$row20[0,18] = 0.9999998074061899
$row20[0,19] = 0.999999931388318
$ws.Range("A20:T20").Value = $row20

$row21 = New-Object 'object[,]' 1,20
$row21[0,0] = 19
$row21[0,1] = "HexGrid-90degTilt5degRes"
$row21[0,2] = 1.000016284539304
$row21[0,3] = 0.9999348587541259
$row21[0,4] = 1.000016284539304
$row21[0,5] = 0.9999348587541259
$row21[0,6] = 1.000016284539304
$row21[0,7] = 1.000016284539304
$row21[0,8] = 1.000043427651147
$row21[0,9] = 0.9999641717090894
$row21[0,10] = 1.000016284539304
$row21[0,11] = 1.000016284539304
$row21[0,12] = 0.999975571646715
$row21[0,13] = 0.999975571646715
$row21[0,14] = 0.9999717716675064
$row21[0,15] = 0.9999891426109114
$row21[0,16] = 0.9999891426109114
$row21[0,17] = 0.9999959280930096
$row21[0,18] = 0.9999959280930096
$row21[0,19] = 0.9999985519553792
$ws.Range("A21:T21").Value = $row21

$row22 = New-Object 'object[,]' 1,20
$row22[0,0] = 20
$row22[0,1] = "HexGrid-90degTilt10degRes"
$row22[0,2] = 1.000057334905025
$row22[0,3] = 0.9997706596311844
$row22[0,4] = 1.000057334905025
$row22[0,5] = 0.9997706596311844
$row22[0,6] = 1.000057334905025
$row22[0,7] = 1.000057334905025
$row22[0,8] = 1.000152895557686
$row22[0,9] = 0.9998738610466495
$row22[0,10] = 1.000057334905025
$row22[0,11] = 1.000057334905025
$row22[0,12] = 0.9999139972681046
$row22[0,13] = 0.9999139972681046
$row22[0,14] = 0.9999006185276196
$row22[0,15] = 0.9999617764804114
$row22[0,16] = 0.9999617764804114
$row22[0,17] = 0.9999856660865647
$row22[0,18] = 0.9999856660865647
$row22[0,19] = 0.9999949034917658
$ws.Range("A22:T22").Value = $row22

$row23 = New-Object 'object[,]' 1,20
$row23[0,0] = 21
$row23[0,1] = "HexGrid-90degTilt15degRes"
$row23[0,2] = 1.000134319885843
$row23[0,3] = 0.9994627189113763
$row23[0,4] = 1.000134319885843
$row23[0,5] = 0.9994627189113763
$row23[0,6] = 1.000134319885843
$row23[0,7] = 1.000134319885843
$row23[0,8] = 1.00035818380463
$row23[0,9] = 0.9997044899509185
$row23[0,10] = 1.000134319885843
$row23[0,11] = 1.000134319885843
$row23[0,12] = 0.9997985193986098
$row23[0,13] = 0.9997985193986098
$row23[0,14] = 0.9997671762493794
$row23[0,15] = 0.9999104528943542
$row23[0,16] = 0.9999104528943542
$row23[0,17] = 0.9999664196422264
$row23[0,18] = 0.9999664196422264
$row23[0,19] = 0.9999880587207425
$ws.Range("A23:T23").Value = $row23

# New rows 20-23 need the same bold/bordered style as the rest of column A
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
